$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Suivi du choix du langage de programmation et de l'IDE
$ws.Range("F5").Value = "Oui"
$ws.Range("G5").Value = "JavaScript"

$ws.Range("F6").Value = "Oui"
$ws.Range("G6").Value = "Visual Studio Code"

# Correction orthographique de l'entete de colonne
$ws.Range("G3").Value = "Information potentielle"

# D16 reprend le meme style (sans gras en trop) que le reste de la colonne
$ws.Range("D16").Font.Bold = $false

# Reformer la formule partagee sur C16:C21 (meme calcul, ecriture groupee)
$ws.Range("C16:C21").Formula = "=IF((D16<>""""),ROW(A1),""-"")"

$ws.Range("G8").Select()
